# Clean up stray apostrophes / ampersands in the organization "name" column (B)
# of the active sheet, per the commit: "cleaned user data, combined user data
# with org ids and names".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    18  = "Ardans - aparna"
    25  = "Bobs super ford - aparna"
    42  = "Erols - aparna"
    46  = "Farrells Ice Cream Parlour - aparna"
    54  = "Grand Hotels  Resorts Ltd - aparna"
    58  = "Hermans World of Sporting Goods - aparna"
    62  = "Hudsons MensWear - aparna"
    78  = "Mervyns - aparna"
    79  = "Mighty Caseys - aparna"
    87  = "Odyssey Records  Tapes - aparna"
    88  = "Oles - aparna"
    99  = "Sandys - aparna"
    107 = "Steves Ice Cream - aparna"
    116 = "Turtles Records  Tapes - aparna"
    117 = "United Oil  Gas Corp. - aparna"
    118 = "United Oil  Gas, Singapore - aparna"
    119 = "United Oil  Gas, UK - aparna"
}

foreach ($row in $replacements.Keys) {
    $ws.Range("B$row").Value = $replacements[$row]
}
